$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the old A:K data to B:L
$ws.Columns.Item(1).Insert()

# --- Header row (row 1) ---
# Old A1 (now B1) was blank; give it the new "segments" header text
$ws.Range("B1").Value2 = "segments"
# Match the header formatting (bold, thin border, centered) used by the other headers
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Segment index column (new column A, rows 2-20) ---
# Each data row gets a 0-based running index alongside its (now shifted-to-B) segment name
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}

# Give column A (rows 2-20) the header-like style (bold, border, centered) that used to live on column B
$ws.Range("B1").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Segment name column (new column B, rows 2-20) loses the old bold/border formatting ---
$ws.Range("C2").Copy()
$ws.Range("B2:B20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "done"
